$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (number format / style) of the last existing data row
# (row 38) down into the new row 39, so the date cell in column A keeps the
# same "YYYY-MM-DD HH:MM:SS" style used throughout column A.
$ws.Range("A38:E38").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new forecast row (winter AVERAGE_10_9 series) with its values.
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.2194017515915414
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.1883185981439661
